# Weekly update: insert 3 new rows of data (week of 2022-02-18, Region
# Metropolitana) above the existing "Angeleno"/"Black Amber" rows, which
# get pushed down from rows 182-191 to rows 185-194.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows starting at row 182, pushing existing data down.
$ws.Rows.Item(182).Resize(3).Insert()

# Common/static column values shared by every data row in this block.
$mercadoId   = 8
$mercado     = "Terminal La Palmera de La Serena"
$region      = "Coquimbo"
$codreg      = 4
$tipo        = "Fruta"
$productoId  = 100103
$producto    = "Frutos de hueso (carozo)"
$categoriaId = 100103002
$categoria   = "Ciruela"
$unidad      = "$/bins (450 kilos)"
$kgUnidad    = 450

function Set-Row {
    param(
        [int]$row,
        [double]$fecha,
        [string]$variedad,
        [string]$calidad,
        [double]$volumen,
        [double]$precioMin,
        [double]$precioMax,
        [double]$precioProm,
        [string]$origen,
        [double]$precioKg
    )

    $ws.Cells.Item($row, 1).Value  = $mercadoId
    $ws.Cells.Item($row, 2).Value  = $mercado
    $ws.Cells.Item($row, 3).Value  = $region
    $ws.Cells.Item($row, 4).Value  = $fecha
    $ws.Cells.Item($row, 5).Value  = $codreg
    $ws.Cells.Item($row, 6).Value  = $tipo
    $ws.Cells.Item($row, 7).Value  = $productoId
    $ws.Cells.Item($row, 8).Value  = $producto
    $ws.Cells.Item($row, 9).Value  = $categoriaId
    $ws.Cells.Item($row, 10).Value = $categoria
    $ws.Cells.Item($row, 11).Value = $variedad
    $ws.Cells.Item($row, 12).Value = $calidad
    $ws.Cells.Item($row, 13).Value = $volumen
    $ws.Cells.Item($row, 14).Value = $precioMin
    $ws.Cells.Item($row, 15).Value = $precioMax
    $ws.Cells.Item($row, 16).Value = $precioProm
    $ws.Cells.Item($row, 17).Value = $unidad
    $ws.Cells.Item($row, 18).Value = $origen
    $ws.Cells.Item($row, 19).Value = $precioKg
    $ws.Cells.Item($row, 20).Value = $kgUnidad
}

Set-Row 182 44610 "Angeleno" "Especial" 10 235000 240000 237500 "Región Metropolitana" 528
Set-Row 183 44610 "Angeleno" "Primera"  10 205000 210000 207500 "Región Metropolitana" 461
Set-Row 184 44610 "Angeleno" "Segunda"  16 175000 180000 177500 "Región Metropolitana" 394
